$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows that only change D (Price) and E (Volume(1h)) values
$ws.Range("D2").Value = "31.272.93"
$ws.Range("E2").Value = "  +3.05%  "

$ws.Range("D3").Value = "1.985.97"
$ws.Range("E3").Value = "  +6.30%  "

$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.63%  "

$ws.Range("D5").Value = "0.7825"
$ws.Range("E5").Value = "  +66.43%  "

$ws.Range("D6").Value = "252.58"
$ws.Range("E6").Value = "  +3.66%  "

$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("D8").Value = "0.3382"
$ws.Range("E8").Value = "  +17.68%  "

$ws.Range("D9").Value = "25.58"
$ws.Range("E9").Value = "  +16.20%  "

$ws.Range("D10").Value = "0.06911"
$ws.Range("E10").Value = "  +7.08%  "

$ws.Range("D11").Value = "0.8317"
$ws.Range("E11").Value = "  +15.48%  "

$ws.Range("D12").Value = "0.08113"
$ws.Range("E12").Value = "  +4.49%  "

$ws.Range("D13").Value = "1.989.80"
$ws.Range("E13").Value = "  +6.32%  "

$ws.Range("D14").Value = "100.50"
$ws.Range("E14").Value = "  +4.90%  "

$ws.Range("D15").Value = "5.429"
$ws.Range("E15").Value = "  +5.95%  "

$ws.Range("D16").Value = "271.12"
$ws.Range("E16").Value = "  -2.75%  "

$ws.Range("D17").Value = "31.283.25"
$ws.Range("E17").Value = "  +3.11%  "

$ws.Range("D18").Value = "13.86"
$ws.Range("E18").Value = "  +6.94%  "

$ws.Range("D19").Value = "0.000007911"
$ws.Range("E19").Value = "  +5.12%  "

$ws.Range("D20").Value = "2.252.82"
$ws.Range("E20").Value = "  +6.50%  "

$ws.Range("D21").Value = "5.701"
$ws.Range("E21").Value = "  +9.24%  "

$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("D24").Value = "6.913"
$ws.Range("E24").Value = "  +11.10%  "

$ws.Range("D25").Value = "9.613"
$ws.Range("E25").Value = "  +6.39%  "

$ws.Range("D26").Value = "164.80"
$ws.Range("E26").Value = "  +0.85%  "

$ws.Range("D27").Value = "0.1455"
$ws.Range("E27").Value = "  +51.72%  "

$ws.Range("D28").Value = "19.68"
$ws.Range("E28").Value = "  +5.60%  "

$ws.Range("D29").Value = "2.168"
$ws.Range("E29").Value = "  +15.71%  "

# Row 30 and 31 swap: PancakeSwap <-> Toncoin
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "1.362"
$ws.Range("E30").Value = "  +3.17%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.565"
$ws.Range("E31").Value = "  +6.64%  "

$ws.Range("D32").Value = "4.546"
$ws.Range("E32").Value = "  +8.23%  "

$ws.Range("D33").Value = "4.308"
$ws.Range("E33").Value = "  +5.17%  "

$ws.Range("D34").Value = "0.05167"
$ws.Range("E34").Value = "  +7.55%  "

$ws.Range("D35").Value = "1.209"
$ws.Range("E35").Value = "  +8.26%  "

$ws.Range("D36").Value = "0.7500"
$ws.Range("E36").Value = "  +8.95%  "

$ws.Range("D37").Value = "2.800"
$ws.Range("E37").Value = "  +3.17%  "

$ws.Range("D38").Value = "1.003"
$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("D39").Value = "0.01992"
$ws.Range("E39").Value = "  +6.28%  "

$ws.Range("D40").Value = "2.921"
$ws.Range("E40").Value = "  +3.93%  "

$ws.Range("D41").Value = "6.611"
$ws.Range("E41").Value = "  +6.51%  "

$ws.Range("D42").Value = "77.96"
$ws.Range("E42").Value = "  +5.19%  "

$ws.Range("D43").Value = "0.4613"
$ws.Range("E43").Value = "  +9.48%  "

$ws.Range("D44").Value = "2.053"
$ws.Range("E44").Value = "  +6.02%  "

$ws.Range("D45").Value = "105.25"
$ws.Range("E45").Value = "  +4.52%  "

$ws.Range("D46").Value = "0.8501"
$ws.Range("E46").Value = "  +3.11%  "

$ws.Range("E47").Value = "  +0.52%  "

$ws.Range("D48").Value = "9.964"
$ws.Range("E48").Value = "  +4.21%  "

$ws.Range("D49").Value = "7.460"
$ws.Range("E49").Value = "  +7.43%  "

# Row 50 and 51 swap: Decentraland <-> Elrond
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "36.34"
$ws.Range("E50").Value = "  +3.17%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.4256"
$ws.Range("E51").Value = "  +8.78%  "
